$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.369.30"
$ws.Range("E2").Value = "  +2.57%  "

$ws.Range("D3").Value = "2.862.06"
$ws.Range("E3").Value = "  +7.88%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'195.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.40%  "

$ws.Range("D6").Value = "'600.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.07%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.15%  "

$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("D10").Value = "2.859.55"
$ws.Range("E10").Value = "  +7.83%  "

$ws.Range("D11").Value = "'0.391"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.62%  "

$ws.Range("E12").Value = "  -1.98%  "

$ws.Range("D13").Value = "'4.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.17%  "

$ws.Range("D14").Value = "3.387.45"
$ws.Range("E14").Value = "  +7.81%  "

$ws.Range("D15").Value = "76.168.07"
$ws.Range("E15").Value = "  +2.56%  "

$ws.Range("D16").Value = "'27.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.91%  "

$ws.Range("E17").Value = "  +2.26%  "

$ws.Range("D18").Value = "2.862.95"
$ws.Range("E18").Value = "  +7.65%  "

$ws.Range("D19").Value = "'9.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").Value = "'12.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.61%  "

$ws.Range("D21").Value = "'382.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.27%  "

$ws.Range("E22").Value = "  +4.18%  "

$ws.Range("E23").Value = "  +2.40%  "

$ws.Range("D24").Value = "'72.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.16%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "3.008.42"
$ws.Range("E26").Value = "  +7.74%  "

$ws.Range("D27").Value = "'4.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.92%  "

$ws.Range("D28").Value = "'9.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.63%  "

$ws.Range("E29").Value = "  +12.59%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("D32").Value = "'515.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "

$ws.Range("D33").Value = "'7.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("E34").Value = "  +5.18%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "'166.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.65%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'20.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.00%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.119"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").Value = "'19.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("D40").Value = "'186.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.47%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "'0.347"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.21%  "

$ws.Range("D43").Value = "'5.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.69%  "

$ws.Range("E44").Value = "  +2.08%  "

$ws.Range("D45").Value = "'1.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.16%  "

$ws.Range("D46").Value = "'0.0902"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.17%  "

$ws.Range("D47").Value = "'40.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("E48").Value = "  +2.72%  "

$ws.Range("D49").Value = "'0.578"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.41%  "

$ws.Range("D50").Value = "'0.668"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.66%  "

$ws.Range("E51").Value = "  +3.89%  "
